$d = $word.ActiveDocument

# The author (Nate) filled in the student number + name on the cover
# page. In the original file this was two separate runs ("Student Name"
# and "s"); replace that whole phrase with a single run containing
# "s5273814 – Nathanael Gazzard" (en dash U+2013 between the two).
$dash = [char]0x2013
$newText = "s5273814 " + $dash + " Nathanael Gazzard"

$d.Content.Find.Execute("Student Names", $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)

Write-Host "Replaced student name/number placeholder."
